$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains two bird/fungi observation records on rows 5 and 7
# whose data was swapped (row 5 should become what row 7 used to contain,
# and vice versa). Columns that already hold identical values on both
# rows (C, D, I, K, T, U, V, W, AD, AE, AG, AT, AW, AX, AY, ...) are left
# untouched on purpose.

# --- Plain value columns (numbers / regular text) -----------------------
# These can be swapped with a straightforward Value2 assignment.
$valueCols = @("A","B","E","F","G","H","P","Q","R","S","AH","AI")

foreach ($col in $valueCols) {
    $addr5 = "$col" + "5"
    $addr7 = "$col" + "7"
    $v5 = $ws.Range($addr5).Value2
    $v7 = $ws.Range($addr7).Value2
    $ws.Range($addr5).Value2 = $v7
    $ws.Range($addr7).Value2 = $v5
}

# --- Date/time columns stored as plain text ------------------------------
# Y, Z, AA, AB hold dates/times formatted as text (e.g. "2023-09-09",
# "12:06"). Assigning such strings straight to Value2 makes Excel
# reinterpret them as real date/time serials, which would not match the
# original (text) storage. Route these through copy / paste-special
# (via a scratch cell well outside the used range) so the values keep
# being stored as text, exactly like the source cells.
$scratch = "AZ100"
$textCols = @("Y","Z","AA","AB")
foreach ($col in $textCols) {
    $addr5 = "$col" + "5"
    $addr7 = "$col" + "7"

    $ws.Range($addr5).Copy()
    $ws.Range($scratch).PasteSpecial()

    $ws.Range($addr7).Copy()
    $ws.Range($addr5).PasteSpecial()

    $ws.Range($scratch).Copy()
    $ws.Range($addr7).PasteSpecial()

    $ws.Range($scratch).Clear()
}

# --- Columns present (empty) only on row 7 --------------------------------
# J7, N7 and AF7 are present but empty on row 7 while the corresponding
# cells on row 5 don't exist at all. After the swap row 7 should no
# longer have them, so clear them out.
$clearFromRow7 = @("J","N","AF")
foreach ($col in $clearFromRow7) {
    $ws.Range("$col" + "7").Clear()
}
